$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab16")

# --- Country footnote swap: South Sudan becomes resource-rich (gains "*"),
# --- Nigeria loses resource-rich status (drops "*") ---
$ws.Range("B34").Value = "South Sudan*"
$ws.Range("B57").Value = "Nigeria"

# Re-shade row 34 (South Sudan) to match the other resource-rich ("*") rows,
# by copying the format from an existing shaded row (row 17, Chad*).
$ws.Range("B17:L17").Copy() | Out-Null
$ws.Range("B34:L34").PasteSpecial(-4122) | Out-Null

# Remove shading from row 57 (Nigeria) to match the other non-resource-rich
# rows, by copying the format from an existing unshaded row (row 33).
$ws.Range("B33:L33").Copy() | Out-Null
$ws.Range("B57:L57").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Recalculated aggregate values (regional/grouping totals impacted by the
# --- resource-rich reclassification of South Sudan and Nigeria) ---
$ws.Range("G66").Value = 86.658456253423793
$ws.Range("C69").Value = 15.8636226446054
$ws.Range("D69").Value = 11.906445762028801
$ws.Range("E69").Value = 1.0677554833562
$ws.Range("F69").Value = 21.138483842768899
$ws.Range("G69").Value = 50.1819099932708
$ws.Range("H69").Value = 47944.4965545417
$ws.Range("I69").Value = 34961.298369259799
$ws.Range("J69").Value = 3215.6203404978401
$ws.Range("K69").Value = 63426.873784666597
$ws.Range("L69").Value = 147365.434117967
$ws.Range("C77").Value = 32.560370557038603
$ws.Range("D77").Value = 21.8740044604709
$ws.Range("E77").Value = 0.0232069675227
$ws.Range("F77").Value = 37.3559679318492
$ws.Range("G77").Value = 75.916260128425094
$ws.Range("H77").Value = 1060330.78059807
$ws.Range("I77").Value = 711763.81499293505
$ws.Range("K77").Value = 1215070.8733763299
$ws.Range("L77").Value = 2450665.4497591201
$ws.Range("K78").Value = 8865173.5920156892
$ws.Range("C80").Value = 31.1740714138123
$ws.Range("D80").Value = 21.532312057950499
$ws.Range("E80").Value = 0.12528186768874
$ws.Range("F80").Value = 32.4100874213567
$ws.Range("G80").Value = 69.233800684390005
$ws.Range("H80").Value = 104091.092051174
$ws.Range("I80").Value = 64169.789951062397
$ws.Range("K80").Value = 109028.32150340499
$ws.Range("L80").Value = 217211.732447407
$ws.Range("C82").Value = 17.158733133683899
$ws.Range("D82").Value = 12.994432623353401
$ws.Range("E82").Value = 0.41938723311391002
$ws.Range("F82").Value = 22.938412179492499
$ws.Range("G82").Value = 67.949785664555904
$ws.Range("H82").Value = 423479.32960454997
$ws.Range("I82").Value = 326857.257629478
$ws.Range("K82").Value = 561852.73076976999
$ws.Range("L82").Value = 1539812.89950952
$ws.Range("C84").Value = 14.8178969842428
$ws.Range("D84").Value = 10.3024188637725
$ws.Range("E84").Value = 1.6486549883613499
$ws.Range("F84").Value = 18.499562318322301
$ws.Range("G84").Value = 76.996827937411794
$ws.Range("H84").Value = 67114.216274803402
$ws.Range("I84").Value = 46199.424906890803
$ws.Range("J84").Value = 7745.80522325114
$ws.Range("K84").Value = 84080.093804070595
$ws.Range("L84").Value = 278449.41508325399
$ws.Range("C86").Value = 17.471919620555301
$ws.Range("D86").Value = 13.223009725918899
$ws.Range("E86").Value = 0.13757303179985
$ws.Range("F86").Value = 23.4687429823408
$ws.Range("G86").Value = 66.096683812586605
$ws.Range("H86").Value = 305541.80790835299
$ws.Range("I86").Value = 232985.155932016
$ws.Range("J86").Value = 2748.7878719267501
$ws.Range("K86").Value = 410834.50838120602
$ws.Range("L86").Value = 1146328.57675356
$ws.Range("C87").Value = 17.405873439022798
$ws.Range("D87").Value = 14.1494564313574
$ws.Range("E87").Value = 0.064273662132409995
$ws.Range("F87").Value = 24.5282441213369
$ws.Range("G87").Value = 65.653652139539403
$ws.Range("H87").Value = 1211871.69673795
$ws.Range("I87").Value = 964190.462912458
$ws.Range("J87").Value = 5360.5409313378796
$ws.Range("K87").Value = 1698738.1704154699
$ws.Range("L87").Value = 4475369.9183563804
$ws.Range("C89").Value = 27.8555550047389
$ws.Range("D89").Value = 16.9741697263566
$ws.Range("E89").Value = 0.025756383972289999
$ws.Range("F89").Value = 32.6498147288751
$ws.Range("G89").Value = 63.1398165360199
$ws.Range("H89").Value = 7339067.6953189503
$ws.Range("I89").Value = 4322548.0306030205
$ws.Range("J89").Value = 4871.8272509420603
$ws.Range("K89").Value = 8691796.6839804593
$ws.Range("L89").Value = 17496631.347262502
$ws.Range("C90").Value = 37.228025090159797
$ws.Range("D90").Value = 22.505844732183601
$ws.Range("E90").Value = 0.17824217458364
$ws.Range("F90").Value = 44.2015619670128
$ws.Range("G90").Value = 112.58471784515
$ws.Range("H90").Value = 21830549.9769127
$ws.Range("I90").Value = 11602192.711858399
$ws.Range("J90").Value = 75516.188972288001
$ws.Range("K90").Value = 26143437.1955228
$ws.Range("L90").Value = 67701175.853523701
$ws.Range("C94").Value = 18.4584244215316
$ws.Range("D94").Value = 14.036845996822199
$ws.Range("E94").Value = 0.38211250438707001
$ws.Range("F94").Value = 19.860868049819199
$ws.Range("G94").Value = 116.432289120996
$ws.Range("H94").Value = 122392.02892256
$ws.Range("I94").Value = 91114.263840997606
$ws.Range("J94").Value = 2780.6067488346498
$ws.Range("K94").Value = 129770.764857951
$ws.Range("L94").Value = 802320.73423742899
$ws.Range("C97").Value = 15.035757654823
$ws.Range("D97").Value = 9.6360384266730694
$ws.Range("E97").Value = 0.62856709420517998
$ws.Range("F97").Value = 18.487876880949798
$ws.Range("G97").Value = 58.569141057820197
$ws.Range("H97").Value = 200047.45941695501
$ws.Range("I97").Value = 129513.726545734
$ws.Range("J97").Value = 8708.4638092624391
$ws.Range("K97").Value = 251042.83271861399
$ws.Range("L97").Value = 709475.29705450195
$ws.Range("C98").Value = 12.944068438567101
$ws.Range("D98").Value = 7.2395879634764402
$ws.Range("E98").Value = 0.17716661137802001
$ws.Range("F98").Value = 17.235649630790601
$ws.Range("G98").Value = 57.1850464737471
$ws.Range("H98").Value = 242676.90723545101
$ws.Range("I98").Value = 121699.74749364
$ws.Range("J98").Value = 3637.3353111986798
$ws.Range("K98").Value = 309550.61364351102
$ws.Range("L98").Value = 953989.13854068296

